$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "I2" = 0.9693716918425304
    "J2" = 0.9693716918425304
    "M2" = 0.4214143333333333
    "N2" = 1.264243
    "O2" = 0.02434128610922473
    "P2" = 0.02434128610922473
    "Q2" = 3.858873912817112
    "R2" = 34.729865215354
    "S2" = 0.02359575369732226
    "T2" = 0.02359575369732226
    "I3" = 0.9693716918425304
    "J3" = 0.9693716918425304
    "N3" = 36.386704
    "O3" = 0.7005766871049885
    "P3" = 0.7005766871049887
    "R3" = 999.5747064061121
    "S3" = 0.6791192084443978
    "T3" = 0.6791192084443979
    "I4" = 0.9693716918425304
    "J4" = 0.9693716918425304
    "M4" = 4.762423333333333
    "N4" = 14.28727
    "O4" = 0.2750820267857866
    "P4" = 0.2750820267857866
    "Q4" = 43.60931679145112
    "R4" = 392.48385112306
    "S4" = 0.2666567297008102
    "T4" = 0.2666567297008102
    "G5" = 0.2893236666666667
    "H5" = 0.867971
    "I5" = 0.03062830815746963
    "J5" = 0.03062830815746962
    "M5" = 0.4214143333333333
    "N5" = 1.264243
    "O5" = 0.02434128610922473
    "P5" = 0.02434128610922473
    "Q5" = 0.1219251401058889
    "R5" = 1.097326260953
    "S5" = 0.0007455324119024701
    "T5" = 0.00074553241190247
    "G6" = 0.2893236666666667
    "H6" = 0.867971
    "I6" = 0.03062830815746963
    "J6" = 0.03062830815746962
    "N6" = 36.386704
    "O6" = 0.7005766871049885
    "P6" = 0.7005766871049887
    "Q6" = 3.509178206398222
    "R6" = 31.582603857584
    "S6" = 0.02145747866059077
    "T6" = 0.02145747866059077
    "G7" = 0.2893236666666667
    "H7" = 0.867971
    "I7" = 0.03062830815746963
    "J7" = 0.03062830815746962
    "M7" = 4.762423333333333
    "N7" = 14.28727
    "O7" = 0.2750820267857866
    "P7" = 0.2750820267857866
    "Q7" = 1.377881781018889
    "R7" = 12.40093602917
    "S7" = 0.008425297084976387
    "T7" = 0.008425297084976387
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
